$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 so the original row 2 data shifts down to row 3.
$ws.Rows.Item(2).Insert()

# New row 2: manager login that passes validation.
$ws.Range("B2").Value = "dypYhev"
$ws.Range("A2").Value = "mngr384654"
$ws.Range("C2").Value = "Passed - Valid Login"

# Former row 2 (now row 3) keeps its UserID/Password but the result changes
# to reflect a failed login (soft assert validation added).
$ws.Range("C3").Value = "Failed - Invalid Login"

# Reflect the selection move that happened while the author was working.
$ws.Range("C6").Select()
